$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '63.159.61'
$ws.Cells.Item(2, 5).Value = '  -1.91%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.127.38'
$ws.Cells.Item(3, 5).Value = '  -0.27%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '595.29'
$ws.Cells.Item(5, 5).Value = '  -2.20%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '136.56'
$ws.Cells.Item(6, 5).Value = '  -4.80%  '
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '3.125.66'
$ws.Cells.Item(8, 5).Value = '  -0.20%  '
$ws.Cells.Item(9, 5).Value = '  -2.46%  '
$ws.Cells.Item(10, 5).Value = '  -3.04%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '5.21'
$ws.Cells.Item(11, 5).Value = '  -3.04%  '
$ws.Cells.Item(12, 5).Value = '  -2.53%  '
$ws.Cells.Item(13, 5).Value = '  -2.80%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '34.27'
$ws.Cells.Item(14, 5).Value = '  -2.99%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.640.32'
$ws.Cells.Item(15, 5).Value = '  -0.39%  '
$ws.Cells.Item(16, 5).Value = '  +2.84%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '63.098.75'
$ws.Cells.Item(17, 5).Value = '  -1.94%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '3.128.45'
$ws.Cells.Item(18, 5).Value = '  -0.43%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.71'
$ws.Cells.Item(19, 5).Value = '  -2.14%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '473.94'
$ws.Cells.Item(20, 5).Value = '  -0.53%  '
$ws.Cells.Item(21, 5).Value = '  -3.48%  '
$ws.Cells.Item(22, 5).Value = '  -2.37%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '7.77'
$ws.Cells.Item(23, 5).Value = '  -0.25%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '86.57'
$ws.Cells.Item(24, 5).Value = '  +0.87%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '13.02'
$ws.Cells.Item(25, 5).Value = '  -3.83%  '
$ws.Cells.Item(26, 5).Value = '  -0.04%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.73'
$ws.Cells.Item(27, 5).Value = '  -1.59%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.16'
$ws.Cells.Item(28, 5).Value = '  -2.39%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.95'
$ws.Cells.Item(29, 5).Value = '  -5.94%  '
$ws.Cells.Item(30, 5).Value = '  -0.25%  '
$ws.Cells.Item(31, 5).Value = '  -0.02%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '26.80'
$ws.Cells.Item(32, 5).Value = '  +0.71%  '
$ws.Cells.Item(33, 5).Value = '  -5.97%  '
$ws.Cells.Item(34, 5).Value = '  -4.12%  '
$ws.Cells.Item(35, 5).Value = '  -2.43%  '
$ws.Cells.Item(36, 5).Value = '  -2.26%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '52.16'
$ws.Cells.Item(37, 5).Value = '  -0.92%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.0₃0707'
$ws.Cells.Item(38, 5).Value = '  -4.51%  '
$ws.Cells.Item(39, 5).Value = '  -1.61%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '422.81'
$ws.Cells.Item(40, 5).Value = '  -6.20%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '8.26'
$ws.Cells.Item(41, 5).Value = '  -0.71%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.70'
$ws.Cells.Item(42, 5).Value = '  -9.50%  '
$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.894.26'
$ws.Cells.Item(43, 5).Value = '  +0.55%  '
$ws.Cells.Item(44, 2).Value = 'Kaspa'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.114'
$ws.Cells.Item(44, 5).Value = '  -3.64%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.264'
$ws.Cells.Item(45, 5).Value = '  +0.77%  '
$ws.Cells.Item(46, 2).Value = 'Fetch.AI'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.14'
$ws.Cells.Item(46, 5).Value = '  -4.14%  '
$ws.Cells.Item(47, 2).Value = 'USDe'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.999'
$ws.Cells.Item(47, 5).Value = '  -0.03%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '25.79'
$ws.Cells.Item(48, 5).Value = '  -2.16%  '
$ws.Cells.Item(49, 5).Value = '  -5.20%  '
$ws.Cells.Item(50, 5).Value = '  -0.89%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '120.06'
$ws.Cells.Item(51, 5).Value = '  -0.56%  '
